# Updated RAD Test Cases and Test Data for Existing Liability.
# Katalon re-ran the "FEINmismatch" test suite; rows 2-16 of the
# FEINmismatch sheet record the Result ("Pass"/"Fail") and run Date/time
# of each test case.  Update them to reflect the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEINmismatch")

# Row -> (Result, Date) for the latest test run.
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Mar 27 15:40:36 EDT 2024"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Wed Mar 27 15:40:43 EDT 2024"

$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Wed Mar 27 14:47:19 EDT 2024"

$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Wed Mar 27 14:47:30 EDT 2024"

$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Wed Mar 27 14:47:40 EDT 2024"

$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Wed Mar 27 14:47:47 EDT 2024"

$ws.Range("A8").Value = "Fail"
$ws.Range("B8").Value = "Wed Mar 27 14:47:56 EDT 2024"

$ws.Range("B9").Value = "Wed Mar 27 14:48:05 EDT 2024"

$ws.Range("B10").Value = "Wed Mar 27 14:48:17 EDT 2024"

$ws.Range("B11").Value = "Wed Mar 27 14:48:28 EDT 2024"

$ws.Range("A12").Value = "Fail"
$ws.Range("B12").Value = "Wed Mar 27 14:48:39 EDT 2024"

$ws.Range("A13").Value = "Fail"
$ws.Range("B13").Value = "Wed Mar 27 14:48:47 EDT 2024"

$ws.Range("A14").Value = "Fail"
$ws.Range("B14").Value = "Wed Mar 27 14:48:55 EDT 2024"

$ws.Range("A15").Value = "Fail"
$ws.Range("B15").Value = "Wed Mar 27 14:49:06 EDT 2024"

$ws.Range("A16").Value = "Fail"
$ws.Range("B16").Value = "Wed Mar 27 14:49:18 EDT 2024"
